$wb = $excel.ActiveWorkbook

# --- Update status text: "Ready for handoff" -> "Handback transform failed" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"

# --- Record the handback transform error detail for each locale ---
$wsZhCn.Range("K3").Value = "Handback file name: nptnpbuv.j10 is different with handoff file name: bb10fd67-ef3f-4be0-b8b4-3c32d2b839fd.dfec1afa3ee45a3457a83c479b71e726437d5576.zh-cn."
$wsDeDe.Range("K3").Value = "Handback file name: nptnpbuv.j10 is different with handoff file name: bb10fd67-ef3f-4be0-b8b4-3c32d2b839fd.dfec1afa3ee45a3457a83c479b71e726437d5576.de-de."
